# Updated symbol list: refresh Price (D) and Volume(1h) (E) columns
# for the rows whose crypto quote changed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "296.17" },
    @{ Cell = "E2"; Value = "-2.28%" },
    @{ Cell = "D3"; Value = "31.12" },
    @{ Cell = "E3"; Value = "-1.93%" },
    @{ Cell = "D4"; Value = "5.114" },
    @{ Cell = "E4"; Value = "-1.90%" },
    @{ Cell = "D5"; Value = "0.07371" },
    @{ Cell = "E5"; Value = "1.41%" },
    @{ Cell = "D6"; Value = "7.698" },
    @{ Cell = "E6"; Value = "-1.22%" },
    @{ Cell = "D7"; Value = "3.741" },
    @{ Cell = "E7"; Value = "-0.10%" },
    @{ Cell = "D8"; Value = "1.642" },
    @{ Cell = "E8"; Value = "12.40%" },
    @{ Cell = "D9"; Value = "0.9175" },
    @{ Cell = "E9"; Value = "1.44%" },
    @{ Cell = "D10"; Value = "0.1675" },
    @{ Cell = "E10"; Value = "0.20%" },
    @{ Cell = "D11"; Value = "0.07111" },
    @{ Cell = "E11"; Value = "-4.65%" },
    @{ Cell = "D12"; Value = "0.08002" },
    @{ Cell = "E12"; Value = "0.78%" },
    @{ Cell = "D13"; Value = "0.02990" },
    @{ Cell = "E13"; Value = "0.48%" },
    @{ Cell = "D14"; Value = "0.09904" },
    @{ Cell = "E14"; Value = "-0.38%" },
    @{ Cell = "D15"; Value = "0.001491" },
    @{ Cell = "E15"; Value = "-0.61%" },
    @{ Cell = "D16"; Value = "0.006161" },
    @{ Cell = "E16"; Value = "-4.59%" },
    @{ Cell = "D17"; Value = "3.449" },
    @{ Cell = "E17"; Value = "-0.57%" },
    @{ Cell = "D18"; Value = "2.229" },
    @{ Cell = "E18"; Value = "0.07%" },
    @{ Cell = "E19"; Value = "-1.85%" },
    @{ Cell = "D20"; Value = "0.1314" },
    @{ Cell = "E20"; Value = "-1.16%" },
    @{ Cell = "D21"; Value = "4.549" },
    @{ Cell = "E21"; Value = "5.84%" },
    @{ Cell = "D22"; Value = "0.04619" },
    @{ Cell = "E22"; Value = "1.80%" },
    @{ Cell = "E23"; Value = "-5.15%" },
    @{ Cell = "D24"; Value = "0.001216" },
    @{ Cell = "E24"; Value = "-0.66%" },
    @{ Cell = "D25"; Value = "0.004426" },
    @{ Cell = "E25"; Value = "0.33%" },
    @{ Cell = "D26"; Value = "0.0001299" },
    @{ Cell = "E26"; Value = "-0.49%" },
    @{ Cell = "D27"; Value = "0.0001874" },
    @{ Cell = "E27"; Value = "7.21%" },
    @{ Cell = "D39"; Value = "0.01683" },
    @{ Cell = "E39"; Value = "1.71%" },
    @{ Cell = "D40"; Value = "0.04402" },
    @{ Cell = "E40"; Value = "-1.85%" },
    @{ Cell = "D41"; Value = "0.007141" },
    @{ Cell = "E41"; Value = "-0.76%" },
    @{ Cell = "D42"; Value = "0.1326" },
    @{ Cell = "E42"; Value = "-1.24%" },
    @{ Cell = "D43"; Value = "0.002138" },
    @{ Cell = "E43"; Value = "-8.61%" },
    @{ Cell = "D44"; Value = "0.01107" },
    @{ Cell = "E44"; Value = "-17.38%" },
    @{ Cell = "D45"; Value = "0.00006015" },
    @{ Cell = "E45"; Value = "-0.95%" },
    @{ Cell = "D46"; Value = "1.857" },
    @{ Cell = "E46"; Value = "-1.86%" },
    @{ Cell = "E47"; Value = "-36.69%" }
)

foreach ($u in $updates) {
    # Force text so numeric-looking strings (and "-x.xx%") keep their
    # exact literal formatting instead of being parsed into a float.
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    # Drop the format override again so no stray style survives on the cell.
    $cell.Style = "Normal"
}

